$wb = $excel.ActiveWorkbook

# --- TestSteps sheet ---
$wsSteps = $wb.Worksheets.Item("TestSteps")
$wsSteps.Activate()

for ($row = 2; $row -le 19; $row++) {
    $wsSteps.Cells.Item($row, 8).Value = "PASS"
}

$wsSteps.Range("H20").Select()

# --- TestCases sheet ---
$wsCases = $wb.Worksheets.Item("TestCases")
$wsCases.Activate()

# Mark the first two test cases as enabled ("Yes") and record a PASS result,
# matching the third test case which already has this set.
$wsCases.Range("C2").Value = "Yes"
$wsCases.Range("D2").Value = "PASS"

$wsCases.Range("C3").Value = "Yes"
$wsCases.Range("D3").Value = "PASS"

$wsCases.Range("A2").Select()
